# Insert a new data row at row 322 (pushing existing rows 322:367 down to 323:368)
# and populate it with the new "Choclo" price record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 322. This shifts rows 322-367 down
# to 323-368, and the (now blank) row 322 inherits formatting from the row
# above it (row 321), matching the style used throughout column D (date fmt).
$ws.Rows.Item(322).Insert()

# Populate the new row 322 with the new record's values.
$ws.Cells.Item(322, 1).Value2  = 8
$ws.Cells.Item(322, 2).Value2  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(322, 3).Value2  = "Coquimbo"
$ws.Cells.Item(322, 4).Value2  = 44474
$ws.Cells.Item(322, 5).Value2  = 4
$ws.Cells.Item(322, 6).Value2  = 100112024
$ws.Cells.Item(322, 7).Value2  = "Choclo"
$ws.Cells.Item(322, 8).Value2  = "Dulce o Americano"
$ws.Cells.Item(322, 9).Value2  = "Primera"
$ws.Cells.Item(322, 10).Value2 = 540
$ws.Cells.Item(322, 11).Value2 = 39000
$ws.Cells.Item(322, 12).Value2 = 40000
$ws.Cells.Item(322, 13).Value2 = 39500
$ws.Cells.Item(322, 14).Value2 = "`$/malla 70 unidades"
$ws.Cells.Item(322, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(322, 16).Value2 = 564
$ws.Cells.Item(322, 17).Value2 = 70
$ws.Cells.Item(322, 18).Value2 = "Hortaliza"
